$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to Text so numeric-looking strings (e.g. '548.76')
# are stored verbatim instead of being auto-converted to numbers by Excel's
# normal numeric-literal detection. Style is reset to Normal afterwards so no
# extra formatting is left behind on the cells.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '61.094.06'
$ws.Range("E2").Value = '  -4.12%  '
$ws.Range("D3").Value = '2.470.99'
$ws.Range("E3").Value = '  -5.32%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '548.76'
$ws.Range("E5").Value = '  -4.33%  '
$ws.Range("D6").Value = '145.98'
$ws.Range("E6").Value = '  -6.41%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '0.601'
$ws.Range("E8").Value = '  -3.68%  '
$ws.Range("D9").Value = '2.469.71'
$ws.Range("E9").Value = '  -5.27%  '
$ws.Range("E10").Value = '  -8.82%  '
$ws.Range("E11").Value = '  -1.73%  '
$ws.Range("E12").Value = '  -7.98%  '
$ws.Range("D13").Value = '0.353'
$ws.Range("E13").Value = '  -7.58%  '
$ws.Range("D14").Value = '26.15'
$ws.Range("E14").Value = '  -6.96%  '
$ws.Range("D15").Value = '2.913.04'
$ws.Range("E15").Value = '  -5.30%  '
$ws.Range("D16").Value = '0.0000164'
$ws.Range("E16").Value = '  -8.76%  '
$ws.Range("D17").Value = '60.927.37'
$ws.Range("E17").Value = '  -4.17%  '
$ws.Range("D18").Value = '2.477.24'
$ws.Range("E18").Value = '  -4.86%  '
$ws.Range("D19").Value = '11.06'
$ws.Range("E19").Value = '  -7.68%  '
$ws.Range("D20").Value = '6.97'
$ws.Range("E20").Value = '  -7.39%  '
$ws.Range("D21").Value = '4.19'
$ws.Range("E21").Value = '  -7.38%  '
$ws.Range("D22").Value = '319.20'
$ws.Range("E22").Value = '  -6.84%  '
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("D24").Value = '63.32'
$ws.Range("E24").Value = '  -5.97%  '
$ws.Range("D25").Value = '1.76'
$ws.Range("E25").Value = '  -3.26%  '
$ws.Range("D26").Value = '0.0₃0992'
$ws.Range("E26").Value = '  -8.15%  '
$ws.Range("D27").Value = '2.576.90'
$ws.Range("E27").Value = '  -5.66%  '
$ws.Range("B28").Value = 'Fetch.AI'
$ws.Range("C28").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D28").Value = '1.49'
$ws.Range("E28").Value = '  -4.68%  '
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").Value = '0.998'
$ws.Range("E29").Value = '  -0.38%  '
$ws.Range("D30").Value = '535.14'
$ws.Range("E30").Value = '  -9.26%  '
$ws.Range("D31").Value = '8.33'
$ws.Range("E31").Value = '  -8.61%  '
$ws.Range("D32").Value = '7.68'
$ws.Range("E32").Value = '  -2.87%  '
$ws.Range("D33").Value = '0.150'
$ws.Range("E33").Value = '  -6.91%  '
$ws.Range("D34").Value = '1.89'
$ws.Range("E34").Value = '  -7.77%  '
$ws.Range("D35").Value = '1.57'
$ws.Range("E35").Value = '  -9.19%  '
$ws.Range("D36").Value = '5.91'
$ws.Range("E36").Value = '  -10.09%  '
$ws.Range("D37").Value = '4.87'
$ws.Range("E37").Value = '  -10.19%  '
$ws.Range("D38").Value = '0.998'
$ws.Range("E38").Value = '  -0.09%  '
$ws.Range("D39").Value = '0.376'
$ws.Range("E39").Value = '  -6.33%  '
$ws.Range("D40").Value = '18.37'
$ws.Range("E40").Value = '  -6.72%  '
$ws.Range("D41").Value = '145.69'
$ws.Range("E41").Value = '  -5.50%  '
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.05%  '
$ws.Range("D43").Value = '1.71'
$ws.Range("E43").Value = '  -8.67%  '
$ws.Range("D44").Value = '39.86'
$ws.Range("E44").Value = '  -3.94%  '
$ws.Range("D45").Value = '2.30'
$ws.Range("E45").Value = '  -7.95%  '
$ws.Range("D46").Value = '147.75'
$ws.Range("E46").Value = '  -5.72%  '
$ws.Range("D47").Value = '3.58'
$ws.Range("E47").Value = '  -7.86%  '
$ws.Range("D48").Value = '20.92'
$ws.Range("E48").Value = '  -11.73%  '
$ws.Range("D49").Value = '0.0532'
$ws.Range("E49").Value = '  -9.67%  '
$ws.Range("D50").Value = '0.585'
$ws.Range("E50").Value = '  -6.95%  '
$ws.Range("E51").Value = '  -5.84%  '

$ws.Range("D2:E51").Style = "Normal"
